# Daily attendance processing - 2025-10-07 08:48:10
# Applies the recorded-attendance refresh to the Session Analysis Results sheet:
#  - reorders "Recorded By" email lists for several sessions
#  - updates aggregate counters / percentages in the side summary tables
#  - flips session B2/#2 (row 103) from Pending to Recorded with its data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-G($row, $text) {
    $ws.Cells.Item($row, 7).Value = $text
}

# --- Row 4 (Year2 A1 HISTOLOGY session 1) ---
Set-G 4 "mariam.noureldin@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"

# --- Row 6: Recorded Sessions counter ---
$ws.Cells.Item(6, 12).Value = 71

# --- Row 7 (Year2 A1 PARASITOLOGY session 1) ---
Set-G 7 "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"

# --- Row 8: Pending Sessions counter ---
$ws.Cells.Item(8, 12).Value = 80

# --- Row 9: Coverage % ---
# (use a leading apostrophe so Excel keeps these as literal text instead of
#  auto-converting the percent-looking string into a numeric percentage)
$ws.Cells.Item(9, 12).Formula = "'46.4%"

# --- Row 10: Average Attendance % ---
$ws.Cells.Item(10, 12).Formula = "'49.3%"

# --- Row 17 (Year2 A1 PHYSIOLOGY session 1) ---
Set-G 17 "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# --- Row 18 (Year2 A1 PHYSIOLOGY session 2) ---
Set-G 18 "Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Row 19: A2 class-statistics side table ---
$ws.Cells.Item(19, 19).Formula = "'47.8%"

# --- Row 20: B2 class-statistics side table ---
$ws.Cells.Item(20, 15).Value = 8
$ws.Cells.Item(20, 17).Value = 9
$ws.Cells.Item(20, 18).Formula = "'47.1%"
$ws.Cells.Item(20, 19).Formula = "'46.8%"

# --- Row 21 (Year2 A2 HISTOLOGY session 1) ---
Set-G 21 "mariam.noureldin@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"

# --- Row 24 ---
Set-G 24 "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"

# --- Row 34 ---
Set-G 34 "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# --- Row 35 ---
Set-G 35 "abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Row 37 ---
Set-G 37 "System, Kerelos.zareef@med.asu.edu.eg"

# --- Row 45 ---
Set-G 45 "mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com"

# --- Row 51 ---
Set-G 51 "naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Row 52 ---
Set-G 52 "yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Row 54 ---
Set-G 54 "System, Kerelos.zareef@med.asu.edu.eg"

# --- Row 62 ---
Set-G 62 "mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com"

# --- Row 68 ---
Set-G 68 "naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Row 69 ---
Set-G 69 "yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Row 72 ---
Set-G 72 "mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 76 ---
Set-G 76 "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

# --- Row 83 ---
Set-G 83 "marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"

# --- Row 85 ---
Set-G 85 "Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

# --- Row 86 ---
Set-G 86 "Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(86, 8).Value = "54/154"

# --- Row 98 ---
Set-G 98 "user@user.com, Walaa.h.ghanima@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, nourhanhosni@med.asu.edu.eg"

# --- Row 99 ---
Set-G 99 "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

# --- Row 102 ---
Set-G 102 "Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

# --- Row 103: session is now recorded (was pending) ---
# Copy the "Recorded" row formatting (green fill, style 4) from the row above
# (same class/subject, previous session) onto row 103, then fill in its data.
$ws.Range("A102:I102").Copy() | Out-Null
$ws.Range("A103:I103").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(103, 7).Value = "ola.m.abdelfattah@med.asu.edu.eg"
$ws.Cells.Item(103, 8).Value = "17/224"
$ws.Cells.Item(103, 9).Value = "Recorded"

# --- Row 109 ---
Set-G 109 "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"

# --- Row 116 ---
Set-G 116 "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 119 ---
Set-G 119 "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg"

# --- Row 126 ---
Set-G 126 "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"

# --- Row 133 ---
Set-G 133 "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 136 ---
Set-G 136 "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg"

# --- Row 143 ---
Set-G 143 "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"

# --- Row 144 ---
Set-G 144 "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"

# --- Row 149 ---
Set-G 149 "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

# --- Row 150 ---
Set-G 150 "marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"

# --- Row 151 ---
Set-G 151 "randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System, nourhanmohamed@med.asu.edu.eg"
